$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bn")
$ws.Range("A155").Value = 44053
Write-Host "Test"
